# AMA format of p values
# Remove the leading zero before the decimal point for p-values and effect
# sizes throughout both result tables (e.g. "0.013" -> ".01", "-0.017" -> ".017"),
# per AMA style guidance.

$wb = $excel.ActiveWorkbook

# --- Table 1 ---
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("F3").Value  = "p = .01"
$ws1.Range("G3").Value  = "V = .34"

$ws1.Range("F4").Value  = "p = .001"
$ws1.Range("G4").Value  = "η² = .16"

$ws1.Range("F5").Value  = "ns (p = .07)"
$ws1.Range("G5").Value  = "V = .27"

$ws1.Range("F6").Value  = "ns (p = .37)"
$ws1.Range("G6").Value  = "V = .17"

$ws1.Range("F7").Value  = "p = .02"
$ws1.Range("G7").Value  = "V = .32"

$ws1.Range("F8").Value  = "p = .02"
$ws1.Range("G8").Value  = "V = .32"

$ws1.Range("F9").Value  = "ns (p = .19)"
$ws1.Range("G9").Value  = "V = .21"

$ws1.Range("F10").Value = "ns (p = .09)"
$ws1.Range("G10").Value = "V = .26"

$ws1.Range("F11").Value = "p = .03"
$ws1.Range("G11").Value = "V = .31"

$ws1.Range("F12").Value = "p = .03"
$ws1.Range("G12").Value = "V = .31"

$ws1.Range("F13").Value = "ns (p = .14)"
$ws1.Range("G13").Value = "V = .23"

$ws1.Range("F14").Value = "ns (p = .31)"
$ws1.Range("G14").Value = "V = .18"

$ws1.Range("F15").Value = "ns (p = .54)"
$ws1.Range("G15").Value = "V = .13"

$ws1.Range("F16").Value = "ns (p = .46)"
$ws1.Range("G16").Value = "V = .14"

$ws1.Range("F17").Value = "ns (p = .23)"
$ws1.Range("G17").Value = "V = .2"

$ws1.Range("F18").Value = "p < .001"
$ws1.Range("G18").Value = "η² = .86"

# --- Table 2 ---
$ws2 = $wb.Worksheets.Item("Table 2")

$ws2.Range("F3").Value  = "ns (p = .67)"
$ws2.Range("G3").Value  = "η² = .017"

$ws2.Range("F4").Value  = "p = .03"
$ws2.Range("G4").Value  = "η² = .071"

$ws2.Range("F5").Value  = "ns (p = .07)"
$ws2.Range("G5").Value  = "η² = .045"

$ws2.Range("F6").Value  = "p = .01"
$ws2.Range("G6").Value  = "η² = .091"

$ws2.Range("F7").Value  = "ns (p = .13)"
$ws2.Range("G7").Value  = "η² = .029"

$ws2.Range("F8").Value  = "ns (p = .07)"
$ws2.Range("G8").Value  = "η² = .045"

$ws2.Range("F9").Value  = "ns (p = .22)"
$ws2.Range("G9").Value  = "η² = .014"

$ws2.Range("F10").Value = "ns (p = .10)"
$ws2.Range("G10").Value = "η² = .037"

$ws2.Range("F11").Value = "p = .01"
$ws2.Range("G11").Value = "η² = .097"

$ws2.Range("F12").Value = "p = .02"
$ws2.Range("G12").Value = "η² = .078"

$ws2.Range("F13").Value = "p = .046"
$ws2.Range("G13").Value = "η² = .058"
